$d = $word.ActiveDocument

function Replace-InParagraph($index, $old, $new) {
    $range = $d.Paragraphs($index).Range
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Paragraph 3 ("Sampling. ...") - citation check update
Replace-InParagraph 3 "Ref-K9Y2M5" "Ref-s405183"
Replace-InParagraph 3 "Ref-G6H1J4" "Ref-s405183"

# Paragraph 4 ("Findings. ... SMI was discovered ...") - citation check update
Replace-InParagraph 4 "Ref-A1B2C3" "Ref-u965124"
Replace-InParagraph 4 "Ref-D4E5F6" "Ref-u965124"

# Paragraph 8 ("Purpose. Breast cancer in younger girls ...") - switched to author-date style
Replace-InParagraph 8 "Ref-A1B2C3" "Smith, 2021"
Replace-InParagraph 8 "Ref-D4E5F6" "Smith, 2021"

# Paragraph 9 ("Data Source. ...") - citation check update
Replace-InParagraph 9 "Ref-J7X2B9" "Ref-u846243"

# Paragraph 10 ("Findings. The results indicated ...") - citation check update
Replace-InParagraph 10 "Ref-J7X2BZ" "Ref-f988337"
Replace-InParagraph 10 "Ref-Q9W3R1" "Ref-f988337"
